$d = $word.ActiveDocument

# --- Edit 1: rewrite opening cutscene paragraph (split into 2 paragraphs + bookmark move) ---
$r1 = $d.Content
$null = $r1.Find.Execute("You are a traveler who is on a trip to Hungary")
$p1 = $r1.Paragraphs(1).Range
$frag1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t xml:space="preserve">You are a traveler who is on a trip to Hungary with a group of friends. It is a hot summer day and your group has been on a non airconditioned bus </w:t></w:r><w:r><w:t xml:space="preserve">all day. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Of course</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the bus isn&#8217;t air conditioned..</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> That&#8217;s what you get for trying to do things on the cheap.</w:t></w:r><w:r><w:t xml:space="preserve"> You are </w:t></w:r><w:r><w:t xml:space="preserve">on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Salgotarjani</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Street, which is known for </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">its </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000A12"/><w:sz w:val="48"/><w:szCs w:val="48"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Jewish</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Cemetery. You just did an audiobook tour of the cemetery, which you have downloaded on your phone, but it&#8217;s all in Hungarian, which you didn&#8217;t realize before you bought it. When the group is done exploring, everyone gets back on the bus, but you </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>have to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> go to the bathroom </w:t></w:r><w:r><w:t xml:space="preserve">because </w:t></w:r><w:r><w:t xml:space="preserve">you were chugging water on the hot, stuffy bus. One of your friends says &#8220;Oh yeah, I think I saw </w:t></w:r><w:r><w:t>a bathroom</w:t></w:r><w:r><w:t xml:space="preserve"> in the cemetery, we&#8217;ll wait for you. But hurry up, Charlie isn&#8217;t feeling well. He always gets sick when he smokes too much, and I tried taking away his lighter earlier today, but he still got sick. He </w:t></w:r><w:r><w:t>has</w:t></w:r><w:r><w:t xml:space="preserve"> medicine with him just in case</w:t></w:r><w:r><w:t xml:space="preserve"> this happens</w:t></w:r><w:r><w:t xml:space="preserve">, but </w:t></w:r><w:r><w:t>apparently it&#8217;s not working for him today</w:t></w:r><w:r><w:t xml:space="preserve">.&#8221; </w:t></w:r><w:r><w:t>Your friend then continues i</w:t></w:r><w:r><w:t xml:space="preserve">n a sarcastically joking voice, &#8220;And be careful&#8230; It&#8217;s getting late and this place might be HAUNTED! </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Haha</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.&#8221; You respond</w:t></w:r><w:r><w:t xml:space="preserve"> dismissively</w:t></w:r><w:r><w:t>, &#8220;Yeah, OKAAAAYYY.&#8221;</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p><w:r><w:t>(this part of the cutscene should probably show the audiobook and its pages to highlight that it is important</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; audio book pics) </w:t></w:r><w:r><w:t>and the empty water bottle is in your inventory by default, to be used to solve the water element puzzle)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1.InsertXML($frag1)

# --- Edit 2: move lastRenderedPageBreak earlier in the timer paragraph ---
$r2 = $d.Content
$null = $r2.Find.Execute("*Once the player leaves the interaction")
$p2 = $r2.Paragraphs(1).Range
$frag2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">*Once the player leaves the interaction, an invisible timer will start, maybe 3 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>minutes.*</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> The player needs to go back to their friends to ask them if they know the word for bathroom. *If the timer runs out, </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>the friends will &#8220;Hey, did you go to the bathroom yet? We&#8217;re ready to leave</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>.&#8221;*</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> When the player interacts with their friends, they ask &#8220;Hey, do you guys know the word for bathroom? The groundskeeper doesn&#8217;t know any English.&#8221; The friends respond, &#8220;We don&#8217;t know much Hungarian, but you&#8217;re in luck with this one, it&#8217;s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>furduszoba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.&#8221; *The text </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>furduszoba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> also appears when the friend says it as a clue that it is important*. *The player now can go back to the groundskeeper and select </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>furduszoba</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from the list of choices.* When the player chooses the correct word, the groundskeeper now points them in the direction of the bathroom, which also has now magically appeared *lightning* in a spot that had nothing there before. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p2.InsertXML($frag2)

Write-Output "done"
